# Defect Log.xlsx - add newly logged defects (rows 5-13) found by
# Nhom Bao, Tuan, Tu and mark the submission note in C2.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Creation date used for every new defect row (2011-10-21 == serial 40837)
$createdDate = 40837

# Row 5 - Media Option_Permission
$ws.Range("B5").Value = $createdDate
$ws.Range("C5").Value = "Media Option_Permission"
$ws.Range("D5").Value = "Phân quyền Media cho user, dư thừa không cần thiết."
$ws.Range("E5").Value = "Error"

# Row 6 - Media Option_Component
$ws.Range("B6").Value = $createdDate
$ws.Range("C6").Value = "Media Option_Component"
$ws.Range("D6").Value = "Cấu hình media, dư thừa không cần thiết."
$ws.Range("E6").Value = "Error"

# Row 7 - Category and section
$ws.Range("B7").Value = $createdDate
$ws.Range("C7").Value = "Category và section"
$ws.Range("D7").Value = "Lỗi giá trị không đồng nhất (giá trị liên kết giữa 2 bảng)"
$ws.Range("E7").Value = "Error"

# Row 8 - Article
$ws.Range("B8").Value = $createdDate
$ws.Range("C8").Value = "Article"
$ws.Range("D8").Value = "Thiếu catid"
$ws.Range("E8").Value = "Error"

# Row 9 - Event trong Article
$ws.Range("B9").Value = $createdDate
$ws.Range("C9").Value = "Event trong Article"
$ws.Range("D9").Value = "Thiếu event select change của combox section"
$ws.Range("E9").Value = "Error"

# Row 10 - Media Manager
$ws.Range("B10").Value = $createdDate
$ws.Range("C10").Value = "Media Manager"
$ws.Range("D10").Value = "chkThumbnailView và chkDetailView không thể dùng checkbox"
$ws.Range("E10").Value = "Error"

# Row 11 - Media Manager
$ws.Range("B11").Value = $createdDate
$ws.Range("C11").Value = "Media Manager"
$ws.Range("D11").Value = "txtFilePath bị dư"
$ws.Range("E11").Value = "Error"

# Row 12 - Media Manager
$ws.Range("B12").Value = $createdDate
$ws.Range("C12").Value = "Media Manager"
$ws.Range("D12").Value = "chkSelectFile không cần thiết"
$ws.Range("E12").Value = "Error"

# Row 13 - Media Manager
$ws.Range("B13").Value = $createdDate
$ws.Range("C13").Value = "Media Manager"
$ws.Range("D13").Value = "Thiếu sự kiện select chọn trong folder"
$ws.Range("E13").Value = "Error"

# C2: submission / group note
$ws.Range("C2").Value = "CMS Click and Change"

# Move the active selection to C2, as left by the editor on save
$ws.Range("C2").Select()
